# Daily attendance processing - 2025-12-18 07:53:41
#
# For every row in the "Recorded By" column (G), whenever the recorded-by
# list starts with "System," (i.e. the automated system listed itself
# first among a comma-separated list of recorders), reverse the order of
# the comma-separated entries so "System" moves to the end of the list
# (case preserved on every token).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val.StartsWith("System,")) {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Length - 1)..0]
        $newVal = $reversed -join ", "
        $cell.Value = $newVal
    }
}
